$p = $ppt.ActivePresentation
try {
    $p.ApplyTheme("Office Theme")
} catch {
    Write-Host "ERR:" $_.Exception.Message
}
